# Rename the inline picture shapes living in the document's headers/footers.
#
# The edit swaps the auto-generated "imageN" display name that Word stores
# on each picture's <wp:docPr>/<pic:cNvPr> pair:
#   - both Pearson-logo pictures (in the "first page" footer and the
#     "primary" footer) go from image1.png -> image2.png
#   - the BTEC logo picture (in the "first page" header) goes from
#     image2.jpg -> image1.jpg
#
# We identify each picture defensively by its AlternativeText (the OOXML
# "descr" attribute) rather than relying on a hard-coded Headers/Footers
# index, since wdHeaderFooterIndex (Primary=1 / FirstPage=2 / EvenPages=3)
# does not line up 1:1 with the physical header1.xml/header2.xml parts.
#
# Renaming has to go through Selection: fetching the shape straight off
# Footers(n)/Headers(n).InlineShapes(i) and assigning .Name on it directly
# raises "addressed block not found" in this host, but selecting the
# shape's Range first and renaming the shape via $word.Selection works.

function Rename-InlineShapeByAltText {
    param(
        [string]$AltText,
        [string]$NewName
    )

    $d = $word.ActiveDocument
    for ($s = 1; $s -le $d.Sections.Count; $s++) {
        $sec = $d.Sections($s)

        for ($hfIndex = 1; $hfIndex -le 3; $hfIndex++) {
            $hdr = $sec.Headers($hfIndex)
            if ($hdr.Exists) {
                for ($i = 1; $i -le $hdr.InlineShapes.Count; $i++) {
                    $shp = $hdr.InlineShapes($i)
                    if ($shp.AlternativeText -eq $AltText) {
                        $shp.Range.Select()
                        $word.Selection.InlineShapes(1).Name = $NewName
                    }
                }
            }

            $ftr = $sec.Footers($hfIndex)
            if ($ftr.Exists) {
                for ($i = 1; $i -le $ftr.InlineShapes.Count; $i++) {
                    $shp = $ftr.InlineShapes($i)
                    if ($shp.AlternativeText -eq $AltText) {
                        $shp.Range.Select()
                        $word.Selection.InlineShapes(1).Name = $NewName
                    }
                }
            }
        }
    }
}

Rename-InlineShapeByAltText -AltText "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" -NewName "image2.png"
Rename-InlineShapeByAltText -AltText "BTec_Logo-Orange" -NewName "image1.jpg"
